# This workbook is a weekly price log (Fruta / Ciruela) kept in reverse-
# chronological-ish order per variety. The edit inserts one new reporting
# entry (two rows: "Especial" and "Primera" quality grades for the new
# "Larry Ann" variety, dated 2022-02-10, origin "Provincia de Curicó") at
# the top of the data block (rows 9-10), pushing all the existing rows
# 9-90 down by two positions (to rows 11-92). The sheet's used range grows
# from A1:T90 to A1:T92 accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at the top of the data block (rows 9 and 10),
# shifting every existing row at/after row 9 down by two.
$ws.Rows("9:10").Insert()

# --- New row 9: Larry Ann / Especial --------------------------------------
$ws.Range("A9").Value = 9
$ws.Range("B9").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C9").Value = "Metropolitana"
$ws.Range("D9").Value = 44602
$ws.Range("E9").Value = 13
$ws.Range("F9").Value = "Fruta"
$ws.Range("G9").Value = 100103
$ws.Range("H9").Value = "Frutos de hueso (carozo)"
$ws.Range("I9").Value = 100103002
$ws.Range("J9").Value = "Ciruela"
$ws.Range("K9").Value = "Larry Ann"
$ws.Range("L9").Value = "Especial"
$ws.Range("M9").Value = 340
$ws.Range("N9").Value = 9000
$ws.Range("O9").Value = 9000
$ws.Range("P9").Value = 9000
$ws.Range("Q9").Value = "`$/caja 15 kilos granel"
$ws.Range("R9").Value = "Provincia de Curicó"
$ws.Range("S9").Value = 600
$ws.Range("T9").Value = 15

# --- New row 10: Larry Ann / Primera --------------------------------------
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C10").Value = "Metropolitana"
$ws.Range("D10").Value = 44602
$ws.Range("E10").Value = 13
$ws.Range("F10").Value = "Fruta"
$ws.Range("G10").Value = 100103
$ws.Range("H10").Value = "Frutos de hueso (carozo)"
$ws.Range("I10").Value = 100103002
$ws.Range("J10").Value = "Ciruela"
$ws.Range("K10").Value = "Larry Ann"
$ws.Range("L10").Value = "Primera"
$ws.Range("M10").Value = 380
$ws.Range("N10").Value = 7500
$ws.Range("O10").Value = 7500
$ws.Range("P10").Value = 7500
$ws.Range("Q10").Value = "`$/caja 15 kilos granel"
$ws.Range("R10").Value = "Provincia de Curicó"
$ws.Range("S10").Value = 500
$ws.Range("T10").Value = 15
